# Final term PPT + Report - Ready for submission
#
# This mirrors the author's last editing pass over the deck:
#  - Slide 2 ("Problem Statement"): the intro sentence is retouched and the
#    redundant word "financial" (right before "news") is removed, which
#    splits the run into three ("...basis of " / "news " / "and social...").
#  - Slide 13 ("Basic Implementation..."): the "SVM" table cell is nudged
#    (re-set) as part of the same proofing pass.
#  - Slide 3 ("Efficient Market Hypothesis") and Slide 9 ("Text
#    Preprocessing") are untouched textually - they were simply clicked
#    through during the final review pass.

$p = $ppt.ActivePresentation

# --- Slide 2: "Financial market analysis on the basis of financial news ..." ---
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange

# Remove the redundant "financial " that precedes "news" (char 43, len 10).
$dup = $tr2.Characters(43, 10)
$dup.Text = ""

# Re-touch "news " as its own run, separating it from the remainder of the
# sentence, matching the author's run split.
$newsRun = $tr2.Characters(43, 5)
$newsRun.Text = "news "

# --- Slide 13: "SVM" table cell in the first results table ---
$s13 = $p.Slides.Item(13)
$tblShp = $s13.Shapes.Item(4)
$tbl = $tblShp.Table
$svmCell = $tbl.Cell(3, 1)
$svmRange = $svmCell.Shape.TextFrame.TextRange
$svmRange.Text = "SVM"

# --- Slide 3: "The efficient market hypothesis ..." (text unchanged) ---
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange
$tr3.Text = $tr3.Text

# --- Slide 9: "Text Preprocessing" title (text unchanged) ---
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(1)
$tr9 = $shp9.TextFrame.TextRange
$tr9.Text = $tr9.Text
